$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 119; this shifts existing rows 119:161 down to 120:162.
$ws.Rows.Item(119).Insert()

# Fill the newly inserted row 119 with its data (copy static columns from the row
# that is now directly below it, row 120, then set the row-specific values).
$ws.Cells.Item(119, 1).Value = 7
$ws.Cells.Item(119, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(119, 3).Value = "Ñuble"
$ws.Cells.Item(119, 4).Value = 44755
$ws.Cells.Item(119, 4).NumberFormat = $ws.Cells.Item(120, 4).NumberFormat
$ws.Cells.Item(119, 5).Value = 16
$ws.Cells.Item(119, 6).Value = 100112045
$ws.Cells.Item(119, 7).Value = "Zapallo"
$ws.Cells.Item(119, 8).Value = "Camote"
$ws.Cells.Item(119, 9).Value = "1a (guarda)"
$ws.Cells.Item(119, 10).Value = 200
$ws.Cells.Item(119, 11).Value = 550
$ws.Cells.Item(119, 12).Value = 600
$ws.Cells.Item(119, 13).Value = 575
$ws.Cells.Item(119, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(119, 15).Value = "Región del Maule"
$ws.Cells.Item(119, 16).Value = 575
$ws.Cells.Item(119, 17).Value = 1
$ws.Cells.Item(119, 18).Value = "Hortaliza"
